$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds decimal-looking strings (e.g. "246.50", "1.00") that must
# stay TEXT, matching the source inlineStr cells. Setting .Value to a plain
# numeric-looking string lets Excel coerce it to a real Number (dropping
# significant trailing zeros, e.g. "246.50" -> 246.5). Prefixing the literal
# with an apostrophe (Excel's standard "treat as text" quote-prefix input)
# keeps it stored as text exactly as typed.

$ws.Range('D2').Value = "'" + '35.924.71'
$ws.Range('E2').Value = '  -1.90%  '

$ws.Range('D3').Value = "'" + '1.993.00'
$ws.Range('E3').Value = '  -2.93%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Value = "'" + '246.50'

$ws.Range('D6').Value = "'" + '0.639'
$ws.Range('E6').Value = '  -3.40%  '

$ws.Range('D7').Value = "'" + '59.69'
$ws.Range('E7').Value = '  +9.09%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').Value = "'" + '58.55'
$ws.Range('E9').Value = '  -3.49%  '

$ws.Range('E10').Value = '  -0.37%  '

$ws.Range('E11').Value = '  -1.63%  '

$ws.Range('E12').Value = '  -2.12%  '

$ws.Range('D13').Value = "'" + '0.947'
$ws.Range('E13').Value = '  -2.34%  '

$ws.Range('D14').Value = "'" + '14.75'
$ws.Range('E14').Value = '  -0.21%  '

$ws.Range('D15').Value = "'" + '2.283.55'
$ws.Range('E15').Value = '  -2.97%  '

$ws.Range('D16').Value = "'" + '5.34'
$ws.Range('E16').Value = '  -2.43%  '

$ws.Range('D17').Value = "'" + '19.35'
$ws.Range('E17').Value = '  +13.02%  '

$ws.Range('D18').Value = "'" + '1.989.10'
$ws.Range('E18').Value = '  -3.42%  '

$ws.Range('D19').Value = "'" + '35.844.45'
$ws.Range('E19').Value = '  -1.95%  '

$ws.Range('D20').Value = "'" + '71.78'
$ws.Range('E20').Value = '  -0.31%  '

$ws.Range('E21').Value = '  -1.04%  '

$ws.Range('E22').Value = '  -0.40%  '

$ws.Range('D23').Value = "'" + '233.49'
$ws.Range('E23').Value = '  -2.01%  '

$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = "'" + '2.64'
$ws.Range('E24').Value = '  +15.93%  '

$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = "'" + '1.00'
$ws.Range('E25').Value = '  -0.02%  '

$ws.Range('E26').Value = '  -4.17%  '

$ws.Range('D27').Value = "'" + '9.67'
$ws.Range('E27').Value = '  +4.83%  '

$ws.Range('D28').Value = "'" + '165.18'
$ws.Range('E28').Value = '  -0.19%  '

$ws.Range('D29').Value = "'" + '19.40'
$ws.Range('E29').Value = '  -3.43%  '

$ws.Range('E30').Value = '  -1.83%  '

$ws.Range('D31').Value = "'" + '4.94'
$ws.Range('E31').Value = '  -2.59%  '

$ws.Range('E32').Value = '  -5.42%  '

$ws.Range('D33').Value = "'" + '0.0979'
$ws.Range('E33').Value = '  +12.33%  '

$ws.Range('E34').Value = '  +1.63%  '

$ws.Range('E35').Value = '  +11.07%  '

$ws.Range('D36').Value = "'" + '4.43'
$ws.Range('E36').Value = '  -1.15%  '

$ws.Range('E37').Value = '  +0.12%  '

$ws.Range('D38').Value = "'" + '1.79'
$ws.Range('E38').Value = '  -1.74%  '

$ws.Range('D39').Value = "'" + '5.74'
$ws.Range('E39').Value = '  +13.76%  '

$ws.Range('E40').Value = '  -1.13%  '

$ws.Range('E41').Value = '  -0.61%  '

$ws.Range('E42').Value = '  -0.62%  '

$ws.Range('D43').Value = "'" + '0.0929'
$ws.Range('E43').Value = '  +1.63%  '

$ws.Range('E44').Value = '  -0.16%  '

$ws.Range('E45').Value = '  +3.79%  '

$ws.Range('D46').Value = "'" + '94.04'
$ws.Range('E46').Value = '  -0.54%  '

$ws.Range('D47').Value = "'" + '7.82'
$ws.Range('E47').Value = '  +4.19%  '

$ws.Range('D48').Value = "'" + '1.368.61'
$ws.Range('E48').Value = '  -3.42%  '

$ws.Range('E49').Value = '  -0.65%  '

$ws.Range('E50').Value = '  +2.79%  '

$ws.Range('D51').Value = "'" + '46.89'
$ws.Range('E51').Value = '  +2.54%  '
